# "expresión regular de matrícula" -- add the missing license-plate
# pattern line ("AA-00-000") to the "Matrículas vehiculares" section,
# and bring the style sheet (Normal font color + the ListLabel
# character styles used by the patrón bullet lists) up to date.

$d = $word.ActiveDocument

# --- 1. Insert the missing pattern text -----------------------------
# The target paragraph is the empty one that sits right before the
# "Estos patrones los encontré..." paragraph (itself right after the
# "000 - AAA" pattern line). Locate it by scanning for that anchor
# text instead of a hard-coded paragraph index, so the script is
# resilient to any small structural differences.
$paraCount = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $cur = $d.Paragraphs.Item($i)
    if ($cur.Range.Text -match "Estos patrones") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 1) {
    $targetPara = $d.Paragraphs.Item($anchorIndex - 1)
    if ($targetPara.Range.Text -match "^\s*$") {
        $targetPara.Range.Text = "AA-00-000"
    }
}

# --- 2. Normal style: explicit font colour instead of "auto" --------
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.Font.Color = 655360

# --- 3. Register the ListLabel character styles (9-26) used by the
#        document's bullet/number lists. ----------------------------
$listLabelStyles = @(
    @{ Name = "ListLabel 9";  Cs = "" },
    @{ Name = "ListLabel 10"; Cs = "Courier New" },
    @{ Name = "ListLabel 11"; Cs = "Wingdings" },
    @{ Name = "ListLabel 12"; Cs = "Symbol" },
    @{ Name = "ListLabel 13"; Cs = "Courier New" },
    @{ Name = "ListLabel 14"; Cs = "Wingdings" },
    @{ Name = "ListLabel 15"; Cs = "Symbol" },
    @{ Name = "ListLabel 16"; Cs = "Courier New" },
    @{ Name = "ListLabel 17"; Cs = "Wingdings" },
    @{ Name = "ListLabel 18"; Cs = "" },
    @{ Name = "ListLabel 19"; Cs = "Courier New" },
    @{ Name = "ListLabel 20"; Cs = "Wingdings" },
    @{ Name = "ListLabel 21"; Cs = "Symbol" },
    @{ Name = "ListLabel 22"; Cs = "Courier New" },
    @{ Name = "ListLabel 23"; Cs = "Wingdings" },
    @{ Name = "ListLabel 24"; Cs = "Symbol" },
    @{ Name = "ListLabel 25"; Cs = "Courier New" },
    @{ Name = "ListLabel 26"; Cs = "Wingdings" }
)

foreach ($def in $listLabelStyles) {
    $newStyle = $d.Styles.Add($def.Name, 2)
    $newStyle.QuickStyle = $true
    $newStyle.Font.NameBi = $def.Cs
}
